$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume symbol data (cryptos.xlsx refresh)
$updates = @{
    'D2' = '303.19'
    'E2' = '5.44%'
    'D3' = '34.86'
    'E3' = '12.53%'
    'E4' = '4.82%'
    'D5' = '0.07807'
    'E5' = '6.68%'
    'D6' = '2.315'
    'E6' = '-2.32%'
    'D7' = '8.011'
    'E7' = '3.79%'
    'D8' = '3.991'
    'E8' = '7.17%'
    'D9' = '0.9286'
    'E9' = '2.62%'
    'E10' = '7.91%'
    'D11' = '0.1830'
    'E11' = '8.39%'
    'D12' = '0.08482'
    'E12' = '3.59%'
    'D13' = '0.03520'
    'E13' = '12.50%'
    'D14' = '0.09902'
    'E14' = '-0.31%'
    'D15' = '0.001479'
    'E15' = '-1.38%'
    'D16' = '0.04596'
    'E16' = '1.85%'
    'D17' = '0.005774'
    'E17' = '0.33%'
    'D19' = '2.110'
    'E19' = '1.46%'
    'D20' = '0.3385'
    'E20' = '1.73%'
    'D21' = '0.1326'
    'E21' = '-0.39%'
    'D22' = '4.541'
    'E22' = '7.74%'
    'E23' = '11.22%'
    'D24' = '0.001221'
    'E24' = '0.96%'
    'D25' = '0.004436'
    'E25' = '6.66%'
    'D26' = '0.0001299'
    'E26' = '-0.13%'
    'E27' = '0.05%'
    'D39' = '0.01765'
    'E39' = '12.23%'
    'D40' = '0.04728'
    'E40' = '6.47%'
    'D41' = '0.007680'
    'E41' = '4.25%'
    'D42' = '0.1407'
    'E42' = '6.11%'
    'D43' = '0.007046'
    'E43' = '-26.09%'
    'D44' = '0.002278'
    'E44' = '1.66%'
    'D45' = '0.009320'
    'E45' = '5.12%'
    'D46' = '0.00005945'
    'E46' = '-2.79%'
    'D47' = '0.00000000749'
    'E47' = '-0.11%'
    'E48' = '5.49%'
    'E49' = '34.70%'
    'E50' = '-0.11%'
    'E51' = '-0.11%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

